# Update the CIBMTR Priority Variables ValueSet workbook (v0.1.6 -> v0.1.7):
#   - bump the Version / Status / Date metadata values
#   - replace the Contact display text and add a distinct second Contact row
#   - add a new Jurisdiction row (pushes Description/Purpose/Copyright/Immutable down one row)
#   - keep the wrap-text / top-vertical alignment formatting used throughout the sheets

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- simple value edits (Property/Value table, row positions 1-11 are stable) ---
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-11-22T12:33:30-06:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- insert a new row for "Jurisdiction" (pushes the remaining rows down by one) ---
$ws.Rows("12:12").Insert()
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

# the inserted row doesn't pick up the shared border/fill formatting automatically,
# so copy it from the row right below (a regular data row)
$ws.Range("A13:B13").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)

# --- re-assert wrap text / top vertical alignment across the whole used range on both sheets ---
$ws.Range("A1:B16").WrapText = $true
$ws.Range("A1:B16").VerticalAlignment = -4160

$ws2 = $wb.Worksheets.Item("Include from LOINC")
$ws2.Range("A1:B281").WrapText = $true
$ws2.Range("A1:B281").VerticalAlignment = -4160
